# Updates cryptos list row data (Price / Volume(1h) columns, plus one
# row re-ranking for Filecoin/OKB) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell as plain text, even when the text looks like a
# number (e.g. "511.04"), so Excel does not auto-convert it to a
# numeric value - these columns are free-form scraped strings.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "60.495.47"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "2.612.39"
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue $ws.Range("D5") "511.04"
$ws.Range("E5").Value = "  -0.60%  "
Set-TextValue $ws.Range("D6") "154.68"
$ws.Range("E6").Value = "  -2.54%  "
Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.14%  "
Set-TextValue $ws.Range("D8") "0.586"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").Value = "2.623.85"
$ws.Range("E9").Value = "  -2.56%  "
Set-TextValue $ws.Range("D10") "6.69"
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("E11").Value = "  -0.99%  "
Set-TextValue $ws.Range("D12") "0.347"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "3.066.85"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "60.439.52"
$ws.Range("E15").Value = "  -0.92%  "
Set-TextValue $ws.Range("D16") "21.65"
$ws.Range("E16").Value = "  -1.46%  "
Set-TextValue $ws.Range("D17") "0.0000141"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "2.611.03"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("E19").Value = "  -1.10%  "
Set-TextValue $ws.Range("D20") "351.56"
$ws.Range("E20").Value = "  +0.21%  "
Set-TextValue $ws.Range("D21") "10.62"
$ws.Range("E21").Value = "  +0.32%  "
Set-TextValue $ws.Range("D22") "6.16"
$ws.Range("E22").Value = "  -1.11%  "
Set-TextValue $ws.Range("D23") "0.999"
$ws.Range("E23").Value = "  -0.07%  "
Set-TextValue $ws.Range("D24") "60.59"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("E26").Value = "  -1.07%  "
Set-TextValue $ws.Range("D27") "0.995"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "0.0₃0845"
$ws.Range("E28").Value = "  -3.47%  "
Set-TextValue $ws.Range("D29") "7.38"
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("E30").Value = "  -0.09%  "
Set-TextValue $ws.Range("D31") "19.47"
$ws.Range("E31").Value = "  -0.87%  "
Set-TextValue $ws.Range("D32") "151.16"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("E33").Value = "  -0.78%  "
Set-TextValue $ws.Range("D34") "5.81"
$ws.Range("E34").Value = "  +0.56%  "
Set-TextValue $ws.Range("D35") "4.01"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("E36").Value = "  -2.71%  "
Set-TextValue $ws.Range("D37") "0.883"
$ws.Range("E37").Value = "  +4.42%  "
$ws.Range("E38").Value = "  -1.81%  "
Set-TextValue $ws.Range("D39") "0.849"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D40") "3.77"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D41") "36.29"
$ws.Range("E41").Value = "  +2.34%  "
Set-TextValue $ws.Range("D42") "294.49"
$ws.Range("E42").Value = "  -6.63%  "
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("E44").Value = "  -0.39%  "
Set-TextValue $ws.Range("D45") "0.996"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  -4.23%  "
Set-TextValue $ws.Range("D47") "19.93"
$ws.Range("E47").Value = "  -1.16%  "
Set-TextValue $ws.Range("D48") "4.89"
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("E49").Value = "  -1.27%  "
Set-TextValue $ws.Range("D50") "10.32"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "2.003.32"
$ws.Range("E51").Value = "  -3.78%  "
